{"js": "// Replace each two-digit-by-two-digit multiplication prompt in the document\n// with its new value. Every prompt text in this document is unique, so a\n// literal (non-regex) body.search() finds exactly the one <w:t> run that\n// needs to change for each pair below (old text -> new text), in document\n// order, matching the commit's OOXML diff.\nconst replacements = [\n  [\"75\u00d783=\", \"62\u00d761=\"],\n  [\"82\u00d716=\", \"54\u00d768=\"],\n  [\"53\u00d785=\", \"64\u00d735=\"],\n  [\"75\u00d756=\", \"70\u00d756=\"],\n  [\"55\u00d746=\", \"29\u00d730=\"],\n  [\"38\u00d764=\", \"98\u00d766=\"],\n  [\"71\u00d720=\", \"45\u00d741=\"],\n  [\"36\u00d746=\", \"41\u00d798=\"],\n  [\"77\u00d794=\", \"22\u00d793=\"],\n  [\"62\u00d727=\", \"81\u00d779=\"],\n  [\"28\u00d761=\", \"71\u00d763=\"],\n  [\"46\u00d751=\", \"99\u00d719=\"],\n  [\"73\u00d775=\", \"25\u00d717=\"],\n  [\"94\u00d724=\", \"83\u00d771=\"],\n  [\"23\u00d787=\", \"44\u00d785=\"],\n  [\"84\u00d722=\", \"25\u00d734=\"],\n  [\"43\u00d725=\", \"90\u00d796=\"],\n  [\"14\u00d781=\", \"47\u00d718=\"],\n  [\"87\u00d789=\", \"43\u00d760=\"],\n  [\"67\u00d775=\", \"77\u00d760=\"],\n  [\"24\u00d731=\", \"53\u00d712=\"],\n  [\"63\u00d767=\", \"74\u00d793=\"],\n  [\"93\u00d727=\", \"98\u00d782=\"],\n  [\"27\u00d748=\", \"72\u00d798=\"],\n  [\"45\u00d744=\", \"80\u00d798=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-by-two-digit multiplication prompt with its new\n# value. Every prompt text in this document is unique, so Find.Execute with\n# MatchCase (no wildcards) locates exactly the one run that needs to change\n# for each pair below (old text -> new text), matching the commit's OOXML\n# diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"75\u00d783=\"; New = \"62\u00d761=\" },\n    @{ Old = \"82\u00d716=\"; New = \"54\u00d768=\" },\n    @{ Old = \"53\u00d785=\"; New = \"64\u00d735=\" },\n    @{ Old = \"75\u00d756=\"; New = \"70\u00d756=\" },\n    @{ Old = \"55\u00d746=\"; New = \"29\u00d730=\" },\n    @{ Old = \"38\u00d764=\"; New = \"98\u00d766=\" },\n    @{ Old = \"71\u00d720=\"; New = \"45\u00d741=\" },\n    @{ Old = \"36\u00d746=\"; New = \"41\u00d798=\" },\n    @{ Old = \"77\u00d794=\"; New = \"22\u00d793=\" },\n    @{ Old = \"62\u00d727=\"; New = \"81\u00d779=\" },\n    @{ Old = \"28\u00d761=\"; New = \"71\u00d763=\" },\n    @{ Old = \"46\u00d751=\"; New = \"99\u00d719=\" },\n    @{ Old = \"73\u00d775=\"; New = \"25\u00d717=\" },\n    @{ Old = \"94\u00d724=\"; New = \"83\u00d771=\" },\n    @{ Old = \"23\u00d787=\"; New = \"44\u00d785=\" },\n    @{ Old = \"84\u00d722=\"; New = \"25\u00d734=\" },\n    @{ Old = \"43\u00d725=\"; New = \"90\u00d796=\" },\n    @{ Old = \"14\u00d781=\"; New = \"47\u00d718=\" },\n    @{ Old = \"87\u00d789=\"; New = \"43\u00d760=\" },\n    @{ Old = \"67\u00d775=\"; New = \"77\u00d760=\" },\n    @{ Old = \"24\u00d731=\"; New = \"53\u00d712=\" },\n    @{ Old = \"63\u00d767=\"; New = \"74\u00d793=\" },\n    @{ Old = \"93\u00d727=\"; New = \"98\u00d782=\" },\n    @{ Old = \"27\u00d748=\"; New = \"72\u00d798=\" },\n    @{ Old = \"45\u00d744=\"; New = \"80\u00d798=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n\n    $found = $find.Execute(\n        [ref]$pair.Old,   # FindText\n        [ref]$true,       # MatchCase\n        [ref]$false,      # MatchWholeWord\n        [ref]$false,      # MatchWildcards\n        [ref]$false,      # MatchSoundsLike\n        [ref]$false,      # MatchAllWordForms\n        [ref]$true,       # Forward\n        [ref]1,           # Wrap (wdFindContinue)\n        [ref]$false,      # Format\n        [ref]$pair.New,   # ReplaceWith\n        [ref]2            # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Search text not found: $($pair.Old)\"\n    }\n}\n"}
